$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 387.5
$ws.Range("I12").Value = 500
$ws.Range("J12").Value = 275
$ws.Range("K12").Value = 500
$ws.Range("L12").Value = 275
$ws.Range("M12").Value = -330
$ws.Range("N12").Value = -615

$ws.Range("H20").Value = 3921
$ws.Range("I20").Value = 3921
$ws.Range("K20").Value = 3921
$ws.Range("M20").Value = -3691

$ws.Range("H31").Value = 57.4
$ws.Range("I31").Value = 57.4
$ws.Range("K31").Value = 172.2
$ws.Range("M31").Value = 57.80000000000001

$ws.Range("H35").Value = 3921
$ws.Range("I35").Value = 3921
$ws.Range("K35").Value = 3921
$ws.Range("M35").Value = -3542

$ws.Range("H38").Value = 281.75
$ws.Range("I38").Value = 281.75
$ws.Range("K38").Value = 845.25
$ws.Range("M38").Value = -473.25

$ws.Range("H39").Value = 354.16666
$ws.Range("I39").Value = 380.7143
$ws.Range("J39").Value = 317
$ws.Range("K39").Value = 1142.1429
$ws.Range("L39").Value = 951
$ws.Range("M39").Value = -846.1428999999998
$ws.Range("N39").Value = -1543

$ws.Range("H94").Value = 5343.5
$ws.Range("I94").Value = 5343.5
$ws.Range("K94").Value = 5343.5
$ws.Range("M94").Value = -4892.5

$ws.Range("H113").Value = 3129.6
$ws.Range("J113").Value = 3816
$ws.Range("L113").Value = 3816
$ws.Range("N113").Value = -10324

$ws.Range("H123").Value = 100000
$ws.Range("I123").Value = 100000
$ws.Range("K123").Value = 100000
$ws.Range("M123").Value = -95100

$ws.Range("H125").Value = 1044
$ws.Range("I125").Value = 1044
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 9396
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -6936
$ws.Range("N125").Value = ""

$ws.Range("H132").Value = 13063.143
$ws.Range("I132").Value = 13063.143
$ws.Range("K132").Value = 39189.429
$ws.Range("M132").Value = -36659.429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 808
$ws.Range("I2").Value = 808
$ws.Range("K2").Value = 808
$ws.Range("M2").Value = -695

$ws.Range("H32").Value = 10441.75
$ws.Range("I32").Value = 7754.727
$ws.Range("K32").Value = 7754.727
$ws.Range("M32").Value = -7467.727

$ws.Range("H61").Value = 17295.166
$ws.Range("I61").Value = 29662.666
$ws.Range("J61").Value = 4927.6665
$ws.Range("K61").Value = 29662.666
$ws.Range("L61").Value = 4927.6665
$ws.Range("M61").Value = -29450.666
$ws.Range("N61").Value = -5351.6665

$ws.Range("H116").Value = 808
$ws.Range("I116").Value = 808
$ws.Range("K116").Value = 808
$ws.Range("M116").Value = 1486

$ws.Range("H136").Value = 17295.166
$ws.Range("I136").Value = 29662.666
$ws.Range("J136").Value = 4927.6665
$ws.Range("K136").Value = 88987.99800000001
$ws.Range("L136").Value = 14782.9995
$ws.Range("M136").Value = -86437.99800000001
$ws.Range("N136").Value = -19882.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 808
$ws.Range("I3").Value = 808
$ws.Range("K3").Value = 808
$ws.Range("M3").Value = -694

$ws.Range("H31").Value = 1000
$ws.Range("J31").Value = 1000
$ws.Range("L31").Value = 1000
$ws.Range("N31").Value = -1504

$ws.Range("H128").Value = 2334.5
$ws.Range("I128").Value = 2334.5
$ws.Range("K128").Value = 7003.5
$ws.Range("M128").Value = -4513.5

$ws.Range("H134").Value = 1840.091
$ws.Range("I134").Value = 1840.091
$ws.Range("K134").Value = 5520.272999999999
$ws.Range("M134").Value = -2985.272999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3176.375
$ws.Range("I16").Value = 902.75
$ws.Range("K16").Value = 902.75
$ws.Range("M16").Value = -615.75

$ws.Range("H31").Value = 2945.1
$ws.Range("I31").Value = 2607.8
$ws.Range("J31").Value = 3057.5334
$ws.Range("K31").Value = 2607.8
$ws.Range("L31").Value = 3057.5334
$ws.Range("M31").Value = -2312.8
$ws.Range("N31").Value = -3647.5334

$ws.Range("H34").Value = 2945.1
$ws.Range("I34").Value = 2607.8
$ws.Range("J34").Value = 3057.5334
$ws.Range("K34").Value = 2607.8
$ws.Range("L34").Value = 3057.5334
$ws.Range("M34").Value = -2405.8
$ws.Range("N34").Value = -3461.5334

$ws.Range("H100").Value = 61666.668
$ws.Range("J100").Value = 61666.668
$ws.Range("L100").Value = 61666.668
$ws.Range("N100").Value = -63830.668

$ws.Range("H113").Value = 3176.375
$ws.Range("I113").Value = 902.75
$ws.Range("K113").Value = 902.75
$ws.Range("M113").Value = 1267.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("I7").Value = 98.5
$ws.Range("J7").Value = 199.75
$ws.Range("K7").Value = 295.5
$ws.Range("L7").Value = 599.25
$ws.Range("M7").Value = -183.5
$ws.Range("N7").Value = -823.25

$ws.Range("H17").Value = 396.2
$ws.Range("I17").Value = 316.66666
$ws.Range("J17").Value = 515.5
$ws.Range("K17").Value = 949.9999799999999
$ws.Range("L17").Value = 1546.5
$ws.Range("M17").Value = -780.9999799999999
$ws.Range("N17").Value = -1884.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 4154.1665
$ws.Range("I22").Value = 666.6667
$ws.Range("K22").Value = 666.6667
$ws.Range("M22").Value = -137.6667

$ws.Range("H132").Value = 3021.3333
$ws.Range("I132").Value = 2494.625
$ws.Range("J132").Value = 4074.75
$ws.Range("K132").Value = 7483.875
$ws.Range("L132").Value = 12224.25
$ws.Range("M132").Value = -4953.875
$ws.Range("N132").Value = -17284.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 4782.6
$ws.Range("I30").Value = 471
$ws.Range("K30").Value = 471
$ws.Range("M30").Value = -363

$ws.Range("H35").Value = 1415.8334
$ws.Range("I35").Value = 998.3333
$ws.Range("K35").Value = 998.3333
$ws.Range("M35").Value = -662.3333

$ws.Range("H40").Value = 2062.25
$ws.Range("I40").Value = 1642.5714
$ws.Range("K40").Value = 1642.5714
$ws.Range("M40").Value = -1506.5714

$ws.Range("H132").Value = 5704
$ws.Range("I132").Value = 4407.75
$ws.Range("K132").Value = 13223.25
$ws.Range("M132").Value = -10693.25

$ws.Range("H136").Value = 58227.57
$ws.Range("I136").Value = 25239
$ws.Range("J136").Value = 102212.336
$ws.Range("K136").Value = 75717
$ws.Range("L136").Value = 306637.008
$ws.Range("M136").Value = -73167
$ws.Range("N136").Value = -311737.008

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1859
$ws.Range("I136").Value = 1859
$ws.Range("K136").Value = 5577
$ws.Range("M136").Value = -3027

